# Fruta / hortaliza, semanal
# Rotates the weekly price-report data across rows 2,3,4,5,7,8,9,10,11,12,13,14
# (row 6 is left untouched) so that each row now shows the figures that, in the
# previous version of the workbook, belonged to a different reporting date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the "before" values for the columns that move (D, L, M, N, O, P, Q, R, S, T)
# for every affected row, before any cell is overwritten.
$cols = @("D","L","M","N","O","P","Q","R","S","T")
$rowsInvolved = @(2,3,4,5,7,8,9,10,11,12,13,14)

$snapshot = @{}
foreach ($r in $rowsInvolved) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowData
}

# Mapping: new row -> old row whose data it should receive
$mapping = @{
    2  = 4
    3  = 10
    4  = 11
    5  = 12
    7  = 2
    8  = 9
    9  = 14
    10 = 5
    11 = 3
    12 = 7
    13 = 8
    14 = 13
}

foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    $src = $snapshot[$oldRow]
    foreach ($c in $cols) {
        $ws.Range("$c$newRow").Value = $src[$c]
    }
}
